$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted before the current row 32,
# pushing the existing rows 32-42 down to 33-43.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly data.
$ws.Range("A32").Value = 2
$ws.Range("B32").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44553
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100112032
$ws.Range("G32").Value = "Zapallo italiano"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 7000
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = 7500
$ws.Range("N32").Value = "$/caja 60 unidades"
$ws.Range("O32").Value = "Provincia de Limarí"
$ws.Range("P32").Value = 125
$ws.Range("Q32").Value = 60
$ws.Range("R32").Value = "Hortaliza"
